$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 25.02000000000047
$ws.Range("H2").Value = 1.580388647153248 / 10000000000000000
$ws.Range("K2").Value = 58.78658903679212
$ws.Range("L2").Value = "[54.3104086417905, 63.26276943179373]"
$ws.Range("O2").Value = 1.628973968528041
$ws.Range("P2").Value = "[1.54092132158058, 1.7170266154755023]"
$ws.Range("S2").Value = 54.97703260455373
$ws.Range("T2").Value = "[52.00994799120679, 57.94411721790066]"
$ws.Range("W2").Value = 18.53333333333368
$ws.Range("X2").Value = 18.18270270270305
$ws.Range("Y2").Value = 18.88396396396432

# Row 3
$ws.Range("E3").Value = 22.87000000000014
$ws.Range("H3").Value = 1.580388647153248 / 10000000000000000
$ws.Range("K3").Value = 59.16543786252052
$ws.Range("L3").Value = "[53.75400569077084, 64.5768700342702]"
$ws.Range("O3").Value = -2.540947811912465
$ws.Range("P3").Value = "[-2.6290004588599274, -2.4528951649650033]"
$ws.Range("S3").Value = 55.57854770479633
$ws.Range("T3").Value = "[52.7478051338232, 58.40929027576946]"
$ws.Range("W3").Value = 9.248728728728784
$ws.Range("X3").Value = 8.92822822822828
$ws.Range("Y3").Value = 9.569229229229288
